# pictures instead of letters
# Replace the "letter" cells (which used special big fonts to render
# kana/roman letters and numbers) with plain-text filenames of the
# pictures that will be used instead, and restyle those cells with the
# normal default font.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: "Most Probable" row
$ws.Range("B2").Value = "Hselect1.jpg"
$ws.Range("C2").Value = "Nselect1.jpg"

# Row 3: "Least Probable" row
$ws.Range("B3").Value = "Hselect1.jpg"
$ws.Range("C3").Value = "Nselect1.jpg"

# The picture-filename cells no longer need the big custom letter fonts;
# switch them back to the workbook's normal/default font.
$letterRange = $ws.Range("B2:C3")
$letterRange.Font.Name = "Calibri"
$letterRange.Font.Size = 11
$letterRange.Font.Bold = $false

# Rows shrink now that they just hold filenames instead of giant letters.
$ws.Rows.Item(2).RowHeight = 37.5
$ws.Rows.Item(3).RowHeight = 37.5

# Restore the printed page setup (paper size / orientation) for the sheet.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Update the remembered selection in the sheet view.
$ws.Range("B7").Select()
